# Move the 7 pending rows from the "New" sheet to the end of the
# "Previously added" sheet (rows 227-233), then clear the "New" sheet's
# data rows (keep only the header row).

$wb  = $excel.ActiveWorkbook
$dst = $wb.Worksheets.Item("Previously added")
$src = $wb.Worksheets.Item("New")

# Row data to migrate (A: link, B: price, C: districtText, D: areaText,
# E: cadastreText, F: date serial). Read straight off the source sheet
# so we copy the exact current content instead of hard-coding it twice.
$startRow = 227

for ($i = 0; $i -le 6; $i++) {
    $srcRow = 2 + $i
    $dstRow = $startRow + $i

    $link = $src.Cells.Item($srcRow, 1).Value2
    $price = $src.Cells.Item($srcRow, 2).Value2
    $district = $src.Cells.Item($srcRow, 3).Value2
    $area = $src.Cells.Item($srcRow, 4).Value2
    $cadastre = $src.Cells.Item($srcRow, 5).Value2
    $date = $src.Cells.Item($srcRow, 6).Value2

    # Write values first. A leading apostrophe on the cadastre number
    # forces text storage even though the content looks numeric.
    $dst.Cells.Item($dstRow, 1).Value = $link
    if ($price -eq $null) { $price = "" }
    $dst.Cells.Item($dstRow, 2).Value = "'" + $price
    $dst.Cells.Item($dstRow, 3).Value = "'" + $district
    $dst.Cells.Item($dstRow, 4).Value = "'" + $area
    if ($cadastre -eq $null) { $cadastre = "" }
    $dst.Cells.Item($dstRow, 5).Value = "'" + $cadastre
    $dst.Cells.Item($dstRow, 6).Value = $date

    # Re-apply the same formatting as the row above it (which already
    # carries the correct styles for link / text / date columns), since
    # assigning .Value can reset a cell's number format.
    $dst.Range("A226:F226").Copy()
    $dst.Range("A" + $dstRow + ":F" + $dstRow).PasteSpecial(-4122)

    # Recreate the hyperlink on column A pointing at the same address.
    $dst.Hyperlinks.Add($dst.Cells.Item($dstRow, 1), $link)

    # Restore formatting once more - adding the hyperlink applies Excel's
    # built-in "Hyperlink" style, which must be overridden back to the
    # sheet's own link style (s=3).
    $dst.Range("A226:F226").Copy()
    $dst.Range("A" + $dstRow + ":F" + $dstRow).PasteSpecial(-4122)
}

# Clear the "New" sheet back down to just its header row: drop the
# hyperlinks first (so no stale relationships survive), then delete the
# now-migrated data rows entirely (not just their contents) so the
# sheet's used range shrinks back to row 1.
$src.Hyperlinks.Delete()
$src.Range("A2:A8").EntireRow.Delete()
